$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 129 (ALC)
$ws.Range("H129").Value = 1114.4419
$ws.Range("J129").Value = 1236.2433
$ws.Range("L129").Value = 3708.7299
$ws.Range("N129").Value = -13708.7299

# Row 132 (ALC)
$ws.Range("H132").Value = 2951.513
$ws.Range("I132").Value = 3023.6333
$ws.Range("J132").Value = 2711.111
$ws.Range("K132").Value = 9070.8999
$ws.Range("L132").Value = 8133.333
$ws.Range("M132").Value = -6540.8999
$ws.Range("N132").Value = -13193.333

# Row 135 (ALC)
$ws.Range("H135").Value = 22729990
$ws.Range("I135").Value = 890.8125
$ws.Range("J135").Value = 83340920
$ws.Range("K135").Value = 8017.3125
$ws.Range("L135").Value = 750068280
$ws.Range("M135").Value = -5482.3125
$ws.Range("N135").Value = -750073350

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 5232.653
$ws.Range("I32").Value = 4122.4316
$ws.Range("J32").Value = 15002.6
$ws.Range("K32").Value = 4122.4316
$ws.Range("L32").Value = 15002.6
$ws.Range("M32").Value = -3835.4316
$ws.Range("N32").Value = -15576.6

# Row 74 (ARM)
$ws.Range("H74").Value = 76926430
$ws.Range("I74").Value = 200002960
$ws.Range("K74").Value = 200002960
$ws.Range("M74").Value = -200002086

# Row 77 (ARM)
$ws.Range("H77").Value = 76926430
$ws.Range("I77").Value = 200002960
$ws.Range("K77").Value = 1000014800
$ws.Range("M77").Value = -1000010432

# Row 132 (ARM)
$ws.Range("H132").Value = 16478.6
$ws.Range("I132").Value = 1978.92
$ws.Range("J132").Value = 52727.8
$ws.Range("K132").Value = 5936.76
$ws.Range("L132").Value = 158183.4
$ws.Range("M132").Value = -3406.76
$ws.Range("N132").Value = -163243.4

$ws = $wb.Worksheets.Item("BSM")
# Row 26 (BSM)
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

# Row 40 (BSM)
$ws.Range("H40").Value = 31474
$ws.Range("J40").Value = 31474
$ws.Range("L40").Value = 31474
$ws.Range("N40").Value = -32004

# Row 80 (BSM)
$ws.Range("H80").Value = 690
$ws.Range("I80").Value = 609.7143
$ws.Range("J80").Value = 727.4666999999999
$ws.Range("K80").Value = 609.7143
$ws.Range("L80").Value = 727.4666999999999
$ws.Range("M80").Value = 388.2857
$ws.Range("N80").Value = -2723.4667

# Row 83 (BSM)
$ws.Range("H83").Value = 690
$ws.Range("I83").Value = 609.7143
$ws.Range("J83").Value = 727.4666999999999
$ws.Range("K83").Value = 3048.5715
$ws.Range("L83").Value = 3637.3335
$ws.Range("M83").Value = 1943.4285
$ws.Range("N83").Value = -13621.3335

$ws = $wb.Worksheets.Item("CRP")
# Row 12 (CRP)
$ws.Range("H12").Value = 995
$ws.Range("I12").Value = 205
$ws.Range("J12").Value = 1785
$ws.Range("K12").Value = 205
$ws.Range("L12").Value = 1785
$ws.Range("M12").Value = -35
$ws.Range("N12").Value = -2125

# Row 31 (CRP)
$ws.Range("H31").Value = 3143.558
$ws.Range("I31").Value = 2868.3076
$ws.Range("J31").Value = 3262.8333
$ws.Range("K31").Value = 2868.3076
$ws.Range("L31").Value = 3262.8333
$ws.Range("M31").Value = -2573.3076
$ws.Range("N31").Value = -3852.8333

# Row 34 (CRP)
$ws.Range("H34").Value = 3143.558
$ws.Range("I34").Value = 2868.3076
$ws.Range("J34").Value = 3262.8333
$ws.Range("K34").Value = 2868.3076
$ws.Range("L34").Value = 3262.8333
$ws.Range("M34").Value = -2666.3076
$ws.Range("N34").Value = -3666.8333

# Row 99 (CRP)
$ws.Range("H99").Value = 35718464
$ws.Range("J99").Value = 100005520
$ws.Range("L99").Value = 100005520
$ws.Range("N99").Value = -100008516

# Row 126 (CRP)
$ws.Range("H126").Value = 35718464
$ws.Range("J126").Value = 100005520
$ws.Range("L126").Value = 300016560
$ws.Range("N126").Value = -300021500

# Row 132 (CRP)
$ws.Range("H132").Value = 3774.2856
$ws.Range("I132").Value = 2848.8572
$ws.Range("K132").Value = 8546.571599999999
$ws.Range("M132").Value = -6016.571599999999

# Row 134 (CRP)
$ws.Range("H134").Value = 1312.875
$ws.Range("I134").Value = 1149.3334
$ws.Range("J134").Value = 1803.5
$ws.Range("K134").Value = 3448.0002
$ws.Range("L134").Value = 5410.5
$ws.Range("M134").Value = -913.0001999999999
$ws.Range("N134").Value = -10480.5

$ws = $wb.Worksheets.Item("CUL")
# Row 23 (CUL)
$ws.Range("H23").Value = 394.33334
$ws.Range("I23").Value = 33.5
$ws.Range("K23").Value = 100.5
$ws.Range("M23").Value = 134.5

# Row 26 (CUL)
$ws.Range("H26").Value = 632.8570999999999
$ws.Range("I26").Value = 632.5
$ws.Range("K26").Value = 1897.5
$ws.Range("M26").Value = -1609.5

# Row 54 (CUL)
$ws.Range("H54").Value = 2573.2144
$ws.Range("I54").Value = 1500
$ws.Range("J54").Value = 4004.1667
$ws.Range("K54").Value = 4500
$ws.Range("L54").Value = 12012.5001
$ws.Range("M54").Value = -3941
$ws.Range("N54").Value = -13130.5001

# Row 70 (CUL)
$ws.Range("H70").Value = 3658.923
$ws.Range("I70").Value = 2812
$ws.Range("K70").Value = 8436
$ws.Range("M70").Value = -8121

# Row 73 (CUL)
$ws.Range("H73").Value = 3658.923
$ws.Range("I73").Value = 2812
$ws.Range("K73").Value = 8436
$ws.Range("M73").Value = -7344

# Row 121 (CUL)
$ws.Range("H121").Value = 937.25
$ws.Range("J121").Value = 1086.625
$ws.Range("L121").Value = 3259.875
$ws.Range("N121").Value = -5879.875

# Row 131 (CUL)
$ws.Range("H131").Value = 723.08
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 723.08
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2169.24
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12249.24

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (GSM)
$ws.Range("H122").Value = 95239790
$ws.Range("I122").Value = 30304874
$ws.Range("K122").Value = 90914622
$ws.Range("M122").Value = -90912172

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (LTW)
$ws.Range("H2").Value = 1230769.2
$ws.Range("I2").Value = 1291666.6
$ws.Range("K2").Value = 1291666.6
$ws.Range("M2").Value = -1291554.6

# Row 12 (LTW)
$ws.Range("H12").Value = 19500000
$ws.Range("J12").Value = 19000000
$ws.Range("L12").Value = 19000000
$ws.Range("N12").Value = -19000340

# Row 22 (LTW)
$ws.Range("H22").Value = 6085.4287
$ws.Range("J22").Value = 5432.8335
$ws.Range("L22").Value = 5432.8335
$ws.Range("N22").Value = -6022.8335

# Row 27 (LTW)
$ws.Range("H27").Value = 6085.4287
$ws.Range("J27").Value = 5432.8335
$ws.Range("L27").Value = 5432.8335
$ws.Range("N27").Value = -5646.8335

# Row 33 (LTW)
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
